$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44312
$ws.Range("J2").Value = 30
$ws.Range("D3").Value = 44243
$ws.Range("J3").Value = 80
$ws.Range("K3").Value = 10000
$ws.Range("L3").Value = 11000
$ws.Range("M3").Value = 10375
$ws.Range("O3").Value = 'Provincia de Quillota'
$ws.Range("P3").Value = 173
$ws.Range("D4").Value = 44200
$ws.Range("J4").Value = 10
$ws.Range("N4").Value = '$/caja 60 unidades'
$ws.Range("O4").Value = 'Provincia de Limarí'
$ws.Range("P4").Value = 150
$ws.Range("Q4").Value = 60
$ws.Range("D5").Value = 44284
$ws.Range("J5").Value = 35
$ws.Range("K5").Value = 10000
$ws.Range("L5").Value = 10000
$ws.Range("M5").Value = 10000
$ws.Range("P5").Value = 167
$ws.Range("D6").Value = 44405
$ws.Range("J6").Value = 45
$ws.Range("K6").Value = 9000
$ws.Range("L6").Value = 9000
$ws.Range("M6").Value = 9000
$ws.Range("N6").Value = '$/caja 50 unidades'
$ws.Range("O6").Value = 'Provincia de Quillota'
$ws.Range("P6").Value = 180
$ws.Range("Q6").Value = 50
$ws.Range("D7").Value = 44291
$ws.Range("J7").Value = 20
$ws.Range("K7").Value = 9000
$ws.Range("L7").Value = 9000
$ws.Range("M7").Value = 9000
$ws.Range("P7").Value = 150
$ws.Range("D8").Value = 44315
$ws.Range("J8").Value = 25
$ws.Range("K8").Value = 10000
$ws.Range("L8").Value = 10000
$ws.Range("M8").Value = 10000
$ws.Range("P8").Value = 167
$ws.Range("D9").Value = 44277
$ws.Range("J9").Value = 25
$ws.Range("K9").Value = 10000
$ws.Range("L9").Value = 10000
$ws.Range("M9").Value = 10000
$ws.Range("P9").Value = 167
$ws.Range("D10").Value = 44333
$ws.Range("J10").Value = 25
$ws.Range("M10").Value = 10400
$ws.Range("O10").Value = 'Provincia de Limarí'
$ws.Range("D11").Value = 44186
$ws.Range("J11").Value = 15
$ws.Range("K11").Value = 7000
$ws.Range("L11").Value = 7000
$ws.Range("M11").Value = 7000
$ws.Range("P11").Value = 117
$ws.Range("D12").Value = 44179
$ws.Range("J12").Value = 15
$ws.Range("K12").Value = 7000
$ws.Range("L12").Value = 7000
$ws.Range("M12").Value = 7000
$ws.Range("P12").Value = 117
